# B6-PowerPoint.pptx maintenance edit
#
# 1) Three tables (on the slides that hold the balance-sheet practice
#    exercises) get switched from the deck's custom default table style
#    to the built-in "Themed Style 2 - Accent 2" table style.
# 2) The presentation's theme colour scheme is re-pointed from the
#    custom "Integral" / "Red Violet" palette to the standard Office
#    theme palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newTableStyleId = "{0B2F3931-C31C-4C79-BBB9-8394B952B9D6}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the theme colour scheme over to the Office palette -----
function HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = HexToRgbInt($officeColors[$i - 1])
}
